# Update countries & provincias Spain
# Refreshes the Covid-19 "Pais" dashboard sheet:
#  - updates the "last updated" timestamp string
#  - refreshes Estados Unidos totals (row 4)
#  - Nigeria's case count overtakes Bolivia's -> rows 107/108 swap countries
#  - Libia's case count grows and overtakes Benin (and others) -> rows 161-169 shift
#    down by one, with Libia's fresh numbers landing on row 161

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 23:22"

# --- Row 4 ---
$ws.Range("B4").Value = 610206
$ws.Range("C4").Value = 23265
$ws.Range("D4").Value = 38520
$ws.Range("E4").Value = 545856
$ws.Range("G4").Value = 2190
$ws.Range("H4").Value = 25830

# --- Row 107 ---
$ws.Range("A107").Value = "Nigeria"
$ws.Range("B107").Value = 362
$ws.Range("C107").Value = 19
$ws.Range("D107").Value = 99
$ws.Range("E107").Value = 252
$ws.Range("F107").Value = 2
$ws.Range("H107").Value = 11

# --- Row 108 ---
$ws.Range("A108").Value = "Bolivia"
$ws.Range("B108").Value = 354
$ws.Range("C108").Value = 24
$ws.Range("D108").Value = 6
$ws.Range("E108").Value = 320
$ws.Range("F108").Value = 3
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 28

# --- Row 161 ---
$ws.Range("A161").Value = "Libia"
$ws.Range("C161").Value = 9
$ws.Range("D161").Value = 9
$ws.Range("E161").Value = 25

# --- Row 162 ---
$ws.Range("A162").Value = "Benin"
$ws.Range("B162").Value = 35
$ws.Range("D162").Value = 18
$ws.Range("E162").Value = 16
$ws.Range("H162").Value = 1

# --- Row 163 ---
$ws.Range("A163").Value = "Eritrea"
$ws.Range("B163").Value = 34
$ws.Range("E163").Value = 34
$ws.Range("H163").Value = 0

# --- Row 164 ---
$ws.Range("A164").Value = "Guam"
$ws.Range("C164").Value = 0
$ws.Range("D164").Value = 0
$ws.Range("E164").Value = 31
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 1

# --- Row 165 ---
$ws.Range("A165").Value = "Sudan"
$ws.Range("C165").Value = 3
$ws.Range("D165").Value = 4
$ws.Range("E165").Value = 23
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 1
$ws.Range("H165").Value = 5

# --- Row 166 ---
$ws.Range("A166").Value = "San Martin (Parte Francesa)"
$ws.Range("B166").Value = 32
$ws.Range("C166").Value = 0
$ws.Range("D166").Value = 11
$ws.Range("E166").Value = 19
$ws.Range("F166").Value = 5
$ws.Range("H166").Value = 2

# --- Row 167 ---
$ws.Range("A167").Value = "Mongolia"
$ws.Range("B167").Value = 30
$ws.Range("C167").Value = 13
$ws.Range("E167").Value = 25
$ws.Range("H167").Value = 0

# --- Row 168 ---
$ws.Range("A168").Value = "Siria"
$ws.Range("B168").Value = 29
$ws.Range("C168").Value = 4
$ws.Range("D168").Value = 5
$ws.Range("E168").Value = 22
$ws.Range("H168").Value = 2

# --- Row 169 ---
$ws.Range("A169").Value = "Mozambique"
$ws.Range("B169").Value = 28
$ws.Range("C169").Value = 7
$ws.Range("D169").Value = 2
$ws.Range("E169").Value = 26
$ws.Range("H169").Value = 0
